# Auto-generated Excel COM-interop script
# Applies updated market/profit values to the Sheets workbook (scheduled runner sync).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 929.5
$ws.Range("I18").Value = 929.5
$ws.Range("K18").Value = 929.5
$ws.Range("M18").Value = -645.5

# Row 33
$ws.Range("H33").Value = 541.6
$ws.Range("I33").Value = 292.42856
$ws.Range("K33").Value = 292.42856
$ws.Range("M33").Value = -63.42856

# Row 43
$ws.Range("H43").Value = 5966
$ws.Range("I43").Value = 4966.6665
$ws.Range("K43").Value = 4966.6665
$ws.Range("M43").Value = -4897.6665

# Row 62
$ws.Range("H62").Value = 600
$ws.Range("I62").Value = 600
$ws.Range("K62").Value = 600
$ws.Range("M62").Value = 24

# Row 65
$ws.Range("H65").Value = 600
$ws.Range("I65").Value = 600
$ws.Range("K65").Value = 3000
$ws.Range("M65").Value = 120

# Row 106
$ws.Range("H106").Value = 1836.6666
$ws.Range("I106").Value = 1836.6666
$ws.Range("K106").Value = 1836.6666
$ws.Range("M106").Value = -1205.6666

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3096.8572
$ws.Range("I2").Value = 1056.4445
$ws.Range("K2").Value = 1056.4445
$ws.Range("M2").Value = -943.4445000000001

# Row 45
$ws.Range("H45").Value = 3107.2307
$ws.Range("I45").Value = 1447
$ws.Range("J45").Value = 3409.0908
$ws.Range("K45").Value = 1447
$ws.Range("L45").Value = 3409.0908
$ws.Range("M45").Value = -1070
$ws.Range("N45").Value = -4163.0908

# Row 103
$ws.Range("H103").Value = 47181
$ws.Range("J103").Value = 47181
$ws.Range("L103").Value = 47181
$ws.Range("N103").Value = -49525

# Row 116
$ws.Range("H116").Value = 3096.8572
$ws.Range("I116").Value = 1056.4445
$ws.Range("K116").Value = 1056.4445
$ws.Range("M116").Value = 1237.5555

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# Row 141
$ws.Range("H141").Value = 80000
$ws.Range("J141").Value = 80000
$ws.Range("L141").Value = 80000
$ws.Range("N141").Value = -90360

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3096.8572
$ws.Range("I3").Value = 1056.4445
$ws.Range("K3").Value = 1056.4445
$ws.Range("M3").Value = -942.4445000000001

# Row 94
$ws.Range("H94").Value = 4833.1113
$ws.Range("J94").Value = 4999.6665
$ws.Range("L94").Value = 4999.6665
$ws.Range("N94").Value = -5901.6665

# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 472.76923
$ws.Range("I22").Value = 494
$ws.Range("K22").Value = 494
$ws.Range("M22").Value = -144

# Row 43
$ws.Range("H43").Value = 25440.75
$ws.Range("J43").Value = 25440.75
$ws.Range("L43").Value = 25440.75
$ws.Range("N43").Value = -25808.75

# Row 99
$ws.Range("H99").Value = 2308.889
$ws.Range("I99").Value = 2700
$ws.Range("K99").Value = 2700
$ws.Range("M99").Value = -1202

# Row 101
$ws.Range("H101").Value = 25440.75
$ws.Range("J101").Value = 25440.75
$ws.Range("L101").Value = 25440.75
$ws.Range("N101").Value = -31930.75

# Row 126
$ws.Range("H126").Value = 2308.889
$ws.Range("I126").Value = 2700
$ws.Range("K126").Value = 8100
$ws.Range("M126").Value = -5630

# Row 141
$ws.Range("H141").Value = 97754.11
$ws.Range("J141").Value = 120965.164
$ws.Range("L141").Value = 120965.164
$ws.Range("N141").Value = -131325.164

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1234
$ws.Range("I5").Value = 1234
$ws.Range("K5").Value = 3702
$ws.Range("M5").Value = -3590

# Row 23
$ws.Range("H23").Value = 327.4
$ws.Range("I23").Value = 100
$ws.Range("J23").Value = 384.25
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 1152.75
$ws.Range("M23").Value = -65
$ws.Range("N23").Value = -1622.75

# Row 26
$ws.Range("H26").Value = 37.333332
$ws.Range("I26").Value = 30
$ws.Range("J26").Value = 52
$ws.Range("K26").Value = 90
$ws.Range("L26").Value = 156
$ws.Range("M26").Value = 198
$ws.Range("N26").Value = -732

# Row 61
$ws.Range("H61").Value = 199
$ws.Range("I61").Value = 199
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 597
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -382
$ws.Range("N61").ClearContents()

# Row 135
$ws.Range("H135").Value = 1234
$ws.Range("I135").Value = 1234
$ws.Range("K135").Value = 11106
$ws.Range("M135").Value = -8571

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2325.8572
$ws.Range("J80").Value = 2546.8333
$ws.Range("L80").Value = 2546.8333
$ws.Range("N80").Value = -4542.8333

# Row 83
$ws.Range("H83").Value = 2325.8572
$ws.Range("J83").Value = 2546.8333
$ws.Range("L83").Value = 12734.1665
$ws.Range("N83").Value = -22718.1665

# Row 101
$ws.Range("H101").Value = 15000
$ws.Range("J101").Value = 15000
$ws.Range("L101").Value = 15000
$ws.Range("N101").Value = -21490

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 592.1667
$ws.Range("I55").Value = 200
$ws.Range("J55").Value = 670.6
$ws.Range("K55").Value = 200
$ws.Range("L55").Value = 670.6
$ws.Range("M55").Value = -27
$ws.Range("N55").Value = -1016.6

# Row 64
$ws.Range("H64").Value = 28074
$ws.Range("J64").Value = 28074
$ws.Range("L64").Value = 28074
$ws.Range("N64").Value = -28524

# Row 67
$ws.Range("H67").Value = 28074
$ws.Range("J67").Value = 28074
$ws.Range("L67").Value = 28074
$ws.Range("N67").Value = -29634

# Row 96
$ws.Range("H96").Value = 35000
$ws.Range("J96").Value = 35000
$ws.Range("L96").Value = 35000
$ws.Range("N96").Value = -40492

# Row 101
$ws.Range("H101").Value = 19696
$ws.Range("J101").Value = 19696
$ws.Range("L101").Value = 19696
$ws.Range("N101").Value = -26186

# Row 132
$ws.Range("H132").Value = 6443.25
$ws.Range("I132").Value = 6443.25
$ws.Range("K132").Value = 19329.75
$ws.Range("M132").Value = -16799.75

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 1268.9166
$ws.Range("I136").Value = 1358.6666
$ws.Range("J136").Value = 999.6667
$ws.Range("K136").Value = 4075.9998
$ws.Range("L136").Value = 2999.0001
$ws.Range("M136").Value = -1525.9998
$ws.Range("N136").Value = -8099.0001

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
